$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @("D2", "38.423.04"),
    @("D3", "2.087.65"),
    @("E3", "  +2.01%  "),
    @("E4", "  +0.02%  "),
    @("D5", "228.18"),
    @("E5", "  -0.17%  "),
    @("D6", "0.612"),
    @("E6", "  +0.65%  "),
    @("D7", "60.74"),
    @("E7", "  +0.11%  "),
    @("E8", "  +0.05%  "),
    @("E9", "  +1.25%  "),
    @("D10", "0.0836"),
    @("E10", "  +1.94%  "),
    @("E11", "  -0.13%  "),
    @("D12", "2.399.16"),
    @("E12", "  +2.09%  "),
    @("D13", "14.83"),
    @("E13", "  +1.31%  "),
    @("D14", "22.37"),
    @("E14", "  +6.47%  "),
    @("D15", "0.785"),
    @("E15", "  +0.48%  "),
    @("D16", "5.43"),
    @("E16", "  +3.87%  "),
    @("D17", "2.092.51"),
    @("E17", "  +1.99%  "),
    @("D18", "38.323.38"),
    @("D19", "71.41"),
    @("E19", "  +2.67%  "),
    @("D20", "6.04"),
    @("E20", "  +2.00%  "),
    @("D21", "0.0₃0833"),
    @("E21", "  +1.14%  "),
    @("D22", "225.46"),
    @("E22", "  +0.56%  "),
    @("E23", "  +0.00%  "),
    @("E24", "  -0.30%  "),
    @("E25", "  +1.46%  "),
    @("D26", "169.69"),
    @("E26", "  +0.91%  "),
    @("D27", "9.43"),
    @("E27", "  +0.97%  "),
    @("D28", "0.135"),
    @("E28", "  +4.95%  "),
    @("D29", "19.03"),
    @("E29", "  +1.12%  "),
    @("E30", "  +8.48%  "),
    @("E31", "  -0.27%  "),
    @("D32", "2.34"),
    @("E32", "  +5.38%  "),
    @("E33", "  +6.80%  "),
    @("D34", "4.50"),
    @("E34", "  +2.30%  "),
    @("D35", "0.0606"),
    @("E35", "  +2.08%  "),
    @("D36", "6.42"),
    @("E36", "  -2.65%  "),
    @("D37", "2.38"),
    @("E37", "  +1.72%  "),
    @("E38", "  +2.28%  "),
    @("E39", "  +0.13%  "),
    @("D40", "18.42"),
    @("E40", "  +0.67%  "),
    @("D41", "1.537.18"),
    @("E41", "  -0.14%  "),
    @("D42", "99.87"),
    @("E42", "  +3.53%  "),
    @("E43", "  +1.37%  "),
    @("E44", "  +1.98%  "),
    @("D45", "2.80"),
    @("E45", "  -0.34%  "),
    @("D46", "7.70"),
    @("E46", "  +8.08%  "),
    @("E47", "  -0.59%  "),
    @("E48", "  +0.68%  "),
    @("E49", "  +2.26%  "),
    @("D50", "2.98"),
    @("E50", "  +1.35%  "),
    @("D51", "2.286.96"),
    @("E51", "  +2.15%  "),
)

foreach ($pair in $changes) {
    $addr = $pair[0]
    $val = $pair[1]
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}